$wb = $excel.ActiveWorkbook

# --- 1. Update the ImportingFilesLink URL text on V_AboutCDSPage!C10 ---
# (Do this BEFORE the DataModelDescription text edit below so that the newly
#  introduced shared string re-uses the slot vacated by the old URL string,
#  keeping shared-string ordering stable for the rest of the workbook.)
$wsAboutCds = $wb.Worksheets.Item("V_AboutCDSPage")
$wsAboutCds.Range("C10").Value = "https://raw.githubusercontent.com/CBIIT/bento-cds-static-content/qa/about/CDSPortalUserGuide.pdf"

# --- 2. Update the DataModelDescription text on V_AboutDataModelPage!B3 ---
# Replace "CRDC Data Hub" with "CRDC Data Submission portal" (both occurrences).
$wsDataModel = $wb.Worksheets.Item("V_AboutDataModelPage")
$newDescription = "The model works with diverse data types, offering robust yet flexible infrastructure while adhering to FAIR (Findability, Accessibility, Interoperability, and Reuse) data principles. CRDC Data Submission portal continues to evolve to meet the data needs for various NCI-funded programs.`n`nAll code necessary to use the Bento platform on which the CRDC Data Submission portal is built is provided in the form of Docker containers. However, the Bento code is also available to the public for research and forking and pull requests. There are a number of CRDC Data Model resources to explore on the CDS GitHub page including the CRDC Core Data Model: https://github.com/CBIIT/cds-model."
$wsDataModel.Range("B3").Value = $newDescription

# --- 3. Update view/selection state on V_AboutCDSPage ---
$wsAboutCds.Range("B15").Select()

# --- 4. Update view/selection state on V_AboutDataModelPage, and leave it as ---
#     the active sheet/tab (matches the workbook's original activeTab).
$wsDataModel.Range("B13").Select()
